$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.204664707183838
$ws.Range("B1").Value = 2.564207315444946
$ws.Range("C1").Value = 3.080844402313232
$ws.Range("D1").Value = 5.448067665100098
$ws.Range("E1").Value = 2.990542888641357
